$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.330.09'
Set-TextValue $ws.Range("D3") '2.712.19'
Set-TextValue $ws.Range("E3") '  +2.38%  '
Set-TextValue $ws.Range("E4") '  -0.07%  '
Set-TextValue $ws.Range("D5") '606.80'
Set-TextValue $ws.Range("E5") '  +1.53%  '
Set-TextValue $ws.Range("D6") '166.72'
Set-TextValue $ws.Range("E6") '  +4.72%  '
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("D8") '0.553'
Set-TextValue $ws.Range("E8") '  +2.59%  '
Set-TextValue $ws.Range("D9") '2.710.94'
Set-TextValue $ws.Range("E9") '  +2.35%  '
Set-TextValue $ws.Range("E10") '  +1.64%  '
Set-TextValue $ws.Range("E11") '  +0.75%  '
Set-TextValue $ws.Range("E12") '  +3.84%  '
Set-TextValue $ws.Range("D14") '28.48'
Set-TextValue $ws.Range("E14") '  +1.64%  '
Set-TextValue $ws.Range("D15") '3.206.33'
Set-TextValue $ws.Range("E15") '  +2.35%  '
Set-TextValue $ws.Range("D17") '68.368.01'
Set-TextValue $ws.Range("E17") '  +0.10%  '
Set-TextValue $ws.Range("D18") '2.714.23'
Set-TextValue $ws.Range("E18") '  +1.78%  '
Set-TextValue $ws.Range("D19") '11.92'
Set-TextValue $ws.Range("E19") '  +3.80%  '
Set-TextValue $ws.Range("D20") '370.85'
Set-TextValue $ws.Range("E20") '  +1.92%  '
Set-TextValue $ws.Range("E21") '  +2.22%  '
Set-TextValue $ws.Range("D22") '4.50'
Set-TextValue $ws.Range("E22") '  +1.86%  '
Set-TextValue $ws.Range("E23") '  +3.78%  '
Set-TextValue $ws.Range("E24") '  +0.41%  '
Set-TextValue $ws.Range("D25") '73.11'
Set-TextValue $ws.Range("E25") '  -1.99%  '
Set-TextValue $ws.Range("D26") '1.00'
Set-TextValue $ws.Range("E26") '  -0.05%  '
Set-TextValue $ws.Range("D27") '10.11'
Set-TextValue $ws.Range("E27") '  +2.55%  '
Set-TextValue $ws.Range("E28") '  +2.27%  '
Set-TextValue $ws.Range("E29") '  +0.64%  '
Set-TextValue $ws.Range("D30") '1.00'
Set-TextValue $ws.Range("E30") '  +0.17%  '
Set-TextValue $ws.Range("D31") '576.87'
Set-TextValue $ws.Range("E31") '  +2.60%  '
Set-TextValue $ws.Range("E32") '  +1.21%  '
Set-TextValue $ws.Range("D33") '1.42'
Set-TextValue $ws.Range("E33") '  +0.92%  '
Set-TextValue $ws.Range("E34") '  +5.84%  '
Set-TextValue $ws.Range("D35") '0.131'
Set-TextValue $ws.Range("E35") '  +1.80%  '
Set-TextValue $ws.Range("E36") '  -0.01%  '
Set-TextValue $ws.Range("E37") '  -3.97%  '
Set-TextValue $ws.Range("D38") '161.51'
Set-TextValue $ws.Range("E38") '  +0.61%  '
Set-TextValue $ws.Range("D39") '19.87'
Set-TextValue $ws.Range("E39") '  +1.07%  '
Set-TextValue $ws.Range("E40") '  +1.80%  '
Set-TextValue $ws.Range("E41") '  -0.10%  '
Set-TextValue $ws.Range("D42") '5.38'
Set-TextValue $ws.Range("E42") '  +0.72%  '
Set-TextValue $ws.Range("D43") '17.99'
Set-TextValue $ws.Range("E44") '  -1.30%  '
Set-TextValue $ws.Range("E45") '  +0.00%  '
Set-TextValue $ws.Range("E46") '  -3.50%  '
Set-TextValue $ws.Range("D47") '40.74'
Set-TextValue $ws.Range("E47") '  +1.18%  '
Set-TextValue $ws.Range("E48") '  +3.70%  '
Set-TextValue $ws.Range("D49") '154.71'
Set-TextValue $ws.Range("E49") '  -2.31%  '
Set-TextValue $ws.Range("D50") '3.90'
Set-TextValue $ws.Range("E50") '  +1.75%  '
Set-TextValue $ws.Range("E51") '  +4.42%  '
